# Add a second column of numeric data (B1:B5) next to the existing
# shared-string column A, then move the active selection down to B6
# (mirrors what Excel leaves as the "current cell" after typing values
# down a column and pressing Enter past the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(111, 222, 333, 444, 555)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B6").Select() | Out-Null
